$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 1329
    "F3" = 1835
    "F4" = 143
    "F6" = 6299
    "F7" = 160
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
